# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (D) and "Correspond Handback DateTime" (G)
# columns for the e22cf724-... row (row 3) on both the zh-cn and de-de sheets, so
# they carry their own handback timestamps instead of sharing the values from the
# bfbc7f52-... row (row 2). The de-de sheet's row 2 values stay as they were.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-27 07:55:40"
$wsZhCn.Range("G3").Value = "2016-01-27 07:56:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-27 07:55:52"
$wsDeDe.Range("G3").Value = "2016-01-27 07:56:54"
